# "update staff info page" - add a "Print: ON/OFF" column (H) to the staff
# info sheet, defaulted to "OFF" for every employee row, and leave the
# selection on the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header cell, formatted like the other header cells (G1 "SEA inform").
$ws.Range("H1").Value = "Print: ON/OFF"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats

# Data cells, formatted like the other boolean/flag cells in column G.
$ws.Range("H2:H4").Value = "OFF"
$ws.Range("G2").Copy()
$ws.Range("H2:H4").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = $false

# Size the new column to fit its contents, same as the rest of the sheet.
$ws.Range("H1:H4").EntireColumn.AutoFit()

# Leave the selection where the author left it after the edit.
[void]$ws.Range("H7").Select()
